# aligner_Error_Code.xlsx — trim the duplicated second line from the two
# "얼라인 설정" / "얼라인 실행" cells, then leave the sheet's selection on D25
# (matching the saved view state in the target workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "얼라인 설정"
$ws.Range("C25").Value = "얼라인 실행"

$ws.Range("D25").Select()
